$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Insert 2 new columns before column D. This shifts existing columns
# D:K (the 8 quarterly data columns) to F:M, and leaves two blank
# columns (D:E) for the two new quarters of data.
# ---------------------------------------------------------------------
$ws.Columns("D:E").Insert()

# ---------------------------------------------------------------------
# The newly inserted D:E columns default to the General style. Copy the
# number/date formatting from the (now shifted) F:G columns into D:E so
# the new columns look like the rest of the data (date format for the
# header row, number format for the data rows). We do this per-block so
# that section-header rows (which only have a single label cell) are
# left untouched, matching the original layout.
# ---------------------------------------------------------------------
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Populate the new D:E columns (two newest quarters) with their data,
# and correct a handful of previously-reported quarterly figures that
# were restated in this update.
# ---------------------------------------------------------------------
# Row 7 (dates)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373

# Row 38 (dates)
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373

# Row 80 (dates)
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373

# Row 8
$ws.Range("D8").Value = 537800
$ws.Range("E8").Value = 518700

# Row 9
$ws.Range("D9").Value = 284500
$ws.Range("E9").Value = 276900

# Row 10
$ws.Range("D10").Value = 253300
$ws.Range("E10").Value = 241800

# Row 12
$ws.Range("D12").Value = 27600
$ws.Range("E12").Value = 40900

# Row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

# Row 14
$ws.Range("D14").Value = 389700
$ws.Range("E14").Value = 0

# Row 15
$ws.Range("D15").Value = 137200
$ws.Range("E15").Value = 143900

# Row 17
$ws.Range("D17").Value = 884500
$ws.Range("E17").Value = 491100
$ws.Range("F17").Value = 511800
$ws.Range("J17").Value = 463800

# Row 18
$ws.Range("D18").Value = -346700
$ws.Range("E18").Value = 27600
$ws.Range("F18").Value = 44500
$ws.Range("J18").Value = 86100

# Row 20
$ws.Range("D20").Value = -40600
$ws.Range("E20").Value = -9200
$ws.Range("F20").Value = -4100
$ws.Range("H20").Value = -19500
$ws.Range("I20").Value = -13100
$ws.Range("J20").Value = -22400

# Row 21
$ws.Range("D21").Value = -250000
$ws.Range("E21").Value = 162300
$ws.Range("H21").Value = 195700
$ws.Range("I21").Value = 225100

# Row 22
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0

# Row 23
$ws.Range("D23").Value = -387300
$ws.Range("E23").Value = 18500
$ws.Range("H23").Value = 66200
$ws.Range("I23").Value = 106700

# Row 24
$ws.Range("D24").Value = 6400
$ws.Range("E24").Value = 1400
$ws.Range("H24").Value = 28700
$ws.Range("I24").Value = 34300

# Row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0

# Row 26
$ws.Range("D26").Value = -393700
$ws.Range("E26").Value = 17100
$ws.Range("H26").Value = 37500
$ws.Range("I26").Value = 72500

# Row 27
$ws.Range("D27").Value = -393700
$ws.Range("E27").Value = 17100
$ws.Range("H27").Value = 37500
$ws.Range("I27").Value = 72500

# Row 28
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0

# Row 29
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0

# Row 30
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0

# Row 31
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0

# Row 32
$ws.Range("D32").Value = 40600
$ws.Range("E32").Value = 9200
$ws.Range("F32").Value = 4100
$ws.Range("H32").Value = 19500
$ws.Range("I32").Value = 13100
$ws.Range("J32").Value = 22400

# Row 33
$ws.Range("D33").Value = -393700
$ws.Range("E33").Value = 17100
$ws.Range("H33").Value = 37500
$ws.Range("I33").Value = 72500

# Row 34
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0

# Row 35
$ws.Range("D35").Value = -393700
$ws.Range("E35").Value = 17100
$ws.Range("H35").Value = 37500
$ws.Range("I35").Value = 72500

# Row 41
$ws.Range("D41").Value = 301800
$ws.Range("E41").Value = 520300

# Row 42
$ws.Range("D42").Value = 82600
$ws.Range("E42").Value = 89300

# Row 43
$ws.Range("D43").Value = 27900
$ws.Range("E43").Value = 45000

# Row 44
$ws.Range("D44").Value = 494200
$ws.Range("E44").Value = 491900

# Row 45
$ws.Range("D45").Value = 166000
$ws.Range("E45").Value = 211200

# Row 46
$ws.Range("D46").Value = 1072500
$ws.Range("E46").Value = 1357700

# Row 47
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0

# Row 48
$ws.Range("D48").Value = 6234300
$ws.Range("E48").Value = 6173000

# Row 49
$ws.Range("D49").Value = 407800
$ws.Range("E49").Value = 696800

# Row 50
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0

# Row 51
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0

# Row 52
$ws.Range("D52").Value = 138300
$ws.Range("E52").Value = 128300

# Row 53
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0

# Row 54
$ws.Range("D54").Value = 7852800
$ws.Range("E54").Value = 8355800

# Row 57
$ws.Range("D57").Value = 310600
$ws.Range("E57").Value = 407200

# Row 58
$ws.Range("D58").Value = 1900
$ws.Range("E58").Value = 2700

# Row 59
$ws.Range("D59").Value = 48900
$ws.Range("E59").Value = 51400

# Row 60
$ws.Range("D60").Value = 361400
$ws.Range("E60").Value = 461300

# Row 61
$ws.Range("D61").Value = 1721300
$ws.Range("E61").Value = 1721500

# Row 62
$ws.Range("D62").Value = 1220100
$ws.Range("E62").Value = 1224100

# Row 63
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0

# Row 64
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0

# Row 65
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0

# Row 66
$ws.Range("D66").Value = 3302800
$ws.Range("E66").Value = 3407000

# Row 68
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0

# Row 69
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0

# Row 70
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0

# Row 71
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0

# Row 72
$ws.Range("D72").Value = -791300
$ws.Range("E72").Value = -374700

# Row 73
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0

# Row 74
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0

# Row 75
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0

# Row 76
$ws.Range("D76").Value = 4550000
$ws.Range("E76").Value = 4948900

# Row 77
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0

# Row 81
$ws.Range("D81").Value = -393700
$ws.Range("E81").Value = 17100
$ws.Range("H81").Value = 37500
$ws.Range("I81").Value = 72500

# Row 83
$ws.Range("D83").Value = 137200
$ws.Range("E83").Value = 143900

# Row 84
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0

# Row 85
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0

# Row 86
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0

# Row 87
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0

# Row 88
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0

# Row 89
$ws.Range("D89").Value = 140300
$ws.Range("E89").Value = 137600

# Row 91
$ws.Range("D91").Value = -342200
$ws.Range("E91").Value = -310600
$ws.Range("I91").Value = -257000
$ws.Range("J91").Value = -192300

# Row 92
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0

# Row 93
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0

# Row 94
$ws.Range("D94").Value = -336400
$ws.Range("E94").Value = -311900

# Row 96
$ws.Range("D96").Value = -20800
$ws.Range("E96").Value = -21100

# Row 97
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0

# Row 98
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0

# Row 99
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0

# Row 100
$ws.Range("D100").Value = -18100
$ws.Range("E100").Value = -14000

# Row 101
$ws.Range("D101").Value = -4200
$ws.Range("E101").Value = 200

# Row 102
$ws.Range("D102").Value = -218400
$ws.Range("E102").Value = -188000
